$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 19, shifting existing rows 19-39 down to 20-40.
$ws.Rows.Item(19).Insert()

# Match the row height used by the other data rows in this table (18.75).
$ws.Rows.Item(19).RowHeight = 18.75

# Populate the new row 19 with the NWWv4 (ammonia water world) data.
$ws.Range("B19").Value2 = "NWWv4"
$ws.Range("C19").Value2 = "WW"
$ws.Range("D19").Value2 = 387
$ws.Range("E19").Value2 = "none"
$ws.Range("F19").Value2 = "ammonia"
$ws.Range("G19").Value2 = "ammonia"
$ws.Range("H19").Value2 = "nitrogen"
$ws.Range("I19").Value2 = "oxigen"
$ws.Range("J19").Value2 = 67.3
$ws.Range("K19").Value2 = 32.700000000000003

# Update the selected cell to match the new state (K19 instead of K18).
$ws.Range("K19").Select()
